# "ajustement staff et mot ministre"
# Update the SOUTIEN sheet: replace the two placeholder names (TBC / TBD)
# with the confirmed names, then leave the selection on C22 (last edited
# cell area) as the workbook was left when saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOUTIEN")

# Write the "Photographers" row name first so the new shared-string table
# keeps the same ordering as the authored workbook (Dominic Lafleur entry
# before the Josyane Bolduc entry).
$ws.Range("C14").Value = "Dominic Lafleur<br/>Justine Boucher"
$ws.Range("C8").Value = "Josyane Bolduc"

# Move the active selection to C22, matching where the author left off.
$ws.Range("C22").Select()
